$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "family"
$ws.Cells.Item(2, 2).Value = "家族|かぞく"
$ws.Cells.Item(3, 1).Value = "grandfather; old man"
$ws.Cells.Item(3, 2).Value = "おじいさん"
$ws.Cells.Item(4, 1).Value = "grandmother; old woman"
$ws.Cells.Item(4, 2).Value = "おばあさん"
$ws.Cells.Item(5, 1).Value = "older brother"
$ws.Cells.Item(5, 2).Value = "お兄さん|おにいさん"
$ws.Cells.Item(6, 1).Value = "older sister"
$ws.Cells.Item(6, 2).Value = "お姉さん|おねえさん"
$ws.Cells.Item(7, 1).Value = "(my) father"
$ws.Cells.Item(7, 2).Value = "父|ちち"
$ws.Cells.Item(8, 1).Value = "(my) mother"
$ws.Cells.Item(8, 2).Value = "母|はは"
$ws.Cells.Item(9, 1).Value = "(my) older brother"
$ws.Cells.Item(9, 2).Value = "兄|あに"
$ws.Cells.Item(10, 1).Value = "(my) older sister"
$ws.Cells.Item(10, 2).Value = "姉|あね"
$ws.Cells.Item(11, 1).Value = "younger brother"
$ws.Cells.Item(11, 2).Value = "弟|おとうと"
$ws.Cells.Item(12, 1).Value = "younger sister"
$ws.Cells.Item(12, 2).Value = "妹|いもうと"
$ws.Cells.Item(13, 1).Value = "brothers and sisters; siblings"
$ws.Cells.Item(13, 2).Value = "兄弟|きょうだい"
$ws.Cells.Item(14, 1).Value = "man"
$ws.Cells.Item(14, 2).Value = "男の人|おとこのひと"
$ws.Cells.Item(15, 1).Value = "woman"
$ws.Cells.Item(15, 2).Value = "女の人|おんなのひと"
$ws.Cells.Item(16, 1).Value = "company"
$ws.Cells.Item(16, 2).Value = "会社|かいしゃ"
$ws.Cells.Item(17, 1).Value = "cafeteria; dining commons"
$ws.Cells.Item(17, 2).Value = "食堂|しょくどう"
$ws.Cells.Item(18, 1).Value = "department store"
$ws.Cells.Item(18, 2).Value = "デパート"
$ws.Cells.Item(19, 1).Value = "hair"
$ws.Cells.Item(19, 2).Value = "髪|かみ"
$ws.Cells.Item(20, 1).Value = "mouth"
$ws.Cells.Item(20, 2).Value = "口|くち"
$ws.Cells.Item(21, 1).Value = "eye"
$ws.Cells.Item(21, 2).Value = "目|め"
$ws.Cells.Item(22, 1).Value = "glasses"
$ws.Cells.Item(22, 2).Value = "眼鏡|めがね"
$ws.Cells.Item(23, 1).Value = "song"
$ws.Cells.Item(23, 2).Value = "歌|うた"
$ws.Cells.Item(24, 1).Value = "club activity"
$ws.Cells.Item(24, 2).Value = "サークル"
$ws.Cells.Item(25, 1).Value = "car"
$ws.Cells.Item(25, 2).Value = "車|くるま"
$ws.Cells.Item(26, 1).Value = "long"
$ws.Cells.Item(26, 2).Value = "長い|ながい"
$ws.Cells.Item(27, 1).Value = "short (length)"
$ws.Cells.Item(27, 2).Value = "短い|みじかい"
$ws.Cells.Item(28, 1).Value = "fast"
$ws.Cells.Item(28, 2).Value = "速い|はやい"
$ws.Cells.Item(29, 1).Value = "tall (stature)"
$ws.Cells.Item(29, 2).Value = "背が高い|せがたかい"
$ws.Cells.Item(30, 1).Value = "short (stature)"
$ws.Cells.Item(30, 2).Value = "背が低い|せがひくい"
$ws.Cells.Item(31, 1).Value = "bright; smart; clever"
$ws.Cells.Item(31, 2).Value = "頭がいい|あたまがいい"
$ws.Cells.Item(32, 1).Value = "cute"
$ws.Cells.Item(32, 2).Value = "かわいい"
$ws.Cells.Item(33, 1).Value = "kind"
$ws.Cells.Item(33, 2).Value = "親切|しんせつ(な)"
$ws.Cells.Item(34, 1).Value = "convenient"
$ws.Cells.Item(34, 2).Value = "便利|べんり(な)"
$ws.Cells.Item(35, 1).Value = "to sing (～を)"
$ws.Cells.Item(35, 2).Value = "歌う|うたう"
$ws.Cells.Item(36, 1).Value = "to put on (a hat) (～を)"
$ws.Cells.Item(36, 2).Value = "かぶる"
$ws.Cells.Item(37, 1).Value = "to put on (items below your waist) (～を)"
$ws.Cells.Item(37, 2).Value = "はく"
$ws.Cells.Item(38, 1).Value = "to get to know (～を)"
$ws.Cells.Item(38, 2).Value = "知る|しる"
$ws.Cells.Item(39, 1).Value = "I know"
$ws.Cells.Item(39, 2).Value = "知っています|しっています"
$ws.Cells.Item(40, 1).Value = "I do not know"
$ws.Cells.Item(40, 2).Value = "知りません|しりません"
$ws.Cells.Item(41, 1).Value = "to live (～に)"
$ws.Cells.Item(41, 2).Value = "住む|すむ"
$ws.Cells.Item(42, 1).Value = "to work"
$ws.Cells.Item(42, 2).Value = "働く|はたらく"
$ws.Cells.Item(43, 1).Value = "to gain weight"
$ws.Cells.Item(43, 2).Value = "太る|ふとる"
$ws.Cells.Item(44, 1).Value = "to be on the heavy side"
$ws.Cells.Item(44, 2).Value = "太っています|ふとっています"
$ws.Cells.Item(45, 1).Value = "to put on (glasses)"
$ws.Cells.Item(45, 2).Value = "(めがねを)かける"
$ws.Cells.Item(46, 1).Value = "to put on (clothes above your waist) (～を)"
$ws.Cells.Item(46, 2).Value = "着る|きる"
$ws.Cells.Item(47, 1).Value = "to lose weight"
$ws.Cells.Item(47, 2).Value = "やせる"
$ws.Cells.Item(48, 1).Value = "to be thin"
$ws.Cells.Item(48, 2).Value = "やせています"
$ws.Cells.Item(49, 1).Value = "to get married (～と)"
$ws.Cells.Item(49, 2).Value = "結婚する|けっこんする"
$ws.Cells.Item(50, 1).Value = "..., but"
$ws.Cells.Item(50, 2).Value = "～が"
$ws.Cells.Item(51, 1).Value = "not...anything"
$ws.Cells.Item(51, 2).Value = "何も|なにも"
$ws.Cells.Item(52, 1).Value = "[counter for people]"
$ws.Cells.Item(52, 2).Value = "～人|～にん"
$ws.Cells.Item(53, 1).Value = "one person"
$ws.Cells.Item(53, 2).Value = "一人|ひとり"
$ws.Cells.Item(54, 1).Value = "two people"
$ws.Cells.Item(54, 2).Value = "二人|ふたり"
$ws.Cells.Item(55, 1).Value = "nothing in particular"
$ws.Cells.Item(55, 2).Value = "別に|べつに"
$ws.Cells.Item(56, 1).Value = "of course"
$ws.Cells.Item(56, 2).Value = "もちろん"
$ws.Cells.Item(57, 1).Value = "if you like"
$ws.Cells.Item(57, 2).Value = "よかったら"
$ws.Cells.Item(58, 1).Value = "how many people"
$ws.Cells.Item(58, 2).Value = "何人|なんにん"
$ws.Cells.Item(59, 1).Value = "one person"
$ws.Cells.Item(59, 2).Value = "一人|ひとり"
$ws.Cells.Item(60, 1).Value = "two people"
$ws.Cells.Item(60, 2).Value = "二人|ふたり"
$ws.Cells.Item(61, 1).Value = "three people"
$ws.Cells.Item(61, 2).Value = "三人|さんにん"
$ws.Cells.Item(62, 1).Value = "four people"
$ws.Cells.Item(62, 2).Value = "四人|よにん"
$ws.Cells.Item(63, 1).Value = "five people"
$ws.Cells.Item(63, 2).Value = "五人|ごにん"
$ws.Cells.Item(64, 1).Value = "six people"
$ws.Cells.Item(64, 2).Value = "六人|ろくにん"
$ws.Cells.Item(65, 1).Value = "seven people"
$ws.Cells.Item(65, 2).Value = "七人|しちにん／ななにん"
$ws.Cells.Item(66, 1).Value = "eight people"
$ws.Cells.Item(66, 2).Value = "八人|はちにん"
$ws.Cells.Item(67, 1).Value = "nine people"
$ws.Cells.Item(67, 2).Value = "九人|きゅうにん"
$ws.Cells.Item(68, 1).Value = "ten people"
$ws.Cells.Item(68, 2).Value = "十人|じゅうにん"
$ws.Cells.Item(69, 1).Value = "hair"
$ws.Cells.Item(69, 2).Value = "髪|かみ"
$ws.Cells.Item(70, 1).Value = "eyebrows"
$ws.Cells.Item(70, 2).Value = "眉毛|まゆげ"
$ws.Cells.Item(71, 1).Value = "ear"
$ws.Cells.Item(71, 2).Value = "耳|みみ"
$ws.Cells.Item(72, 1).Value = "mouth"
$ws.Cells.Item(72, 2).Value = "口|くち"
$ws.Cells.Item(73, 1).Value = "hand"
$ws.Cells.Item(73, 2).Value = "手|て"
$ws.Cells.Item(74, 1).Value = "neck"
$ws.Cells.Item(74, 2).Value = "首|くび"
$ws.Cells.Item(75, 1).Value = "tooth"
$ws.Cells.Item(75, 2).Value = "歯|は"
$ws.Cells.Item(76, 1).Value = "finger"
$ws.Cells.Item(76, 2).Value = "指|ゆび"
$ws.Cells.Item(77, 1).Value = "nose"
$ws.Cells.Item(77, 2).Value = "鼻|はな"
$ws.Cells.Item(78, 1).Value = "eye"
$ws.Cells.Item(78, 2).Value = "目|め"
$ws.Cells.Item(79, 1).Value = "face"
$ws.Cells.Item(79, 2).Value = "顔|かお"
$ws.Cells.Item(80, 1).Value = "head"
$ws.Cells.Item(80, 2).Value = "頭|あたま"
$ws.Cells.Item(81, 1).Value = "shoulder"
$ws.Cells.Item(81, 2).Value = "肩|かた"
$ws.Cells.Item(82, 1).Value = "chest"
$ws.Cells.Item(82, 2).Value = "胸|むね"
$ws.Cells.Item(83, 1).Value = "back (of body)"
$ws.Cells.Item(83, 2).Value = "背中|せなか"
$ws.Cells.Item(84, 1).Value = "belly; stomach"
$ws.Cells.Item(84, 2).Value = "おなか"
$ws.Cells.Item(85, 1).Value = "bottom; buttocks"
$ws.Cells.Item(85, 2).Value = "おしり"
$ws.Cells.Item(86, 1).Value = "foot"
$ws.Cells.Item(86, 2).Value = "足|あし"
$ws.Cells.Item(87, 1).Value = "Tokyo"
$ws.Cells.Item(87, 2).Value = "東京|とうきょう"
$ws.Cells.Item(88, 1).Value = "Kyoko"
$ws.Cells.Item(88, 2).Value = "京子|きょうこ"
$ws.Cells.Item(89, 1).Value = "Kyoto"
$ws.Cells.Item(89, 2).Value = "京都|きょうと"
$ws.Cells.Item(90, 1).Value = "to go to the capital"
$ws.Cells.Item(90, 2).Value = "上京する|じょうきょうする"
$ws.Cells.Item(91, 1).Value = "child"
$ws.Cells.Item(91, 2).Value = "子供|こども"
$ws.Cells.Item(92, 1).Value = "girl"
$ws.Cells.Item(92, 2).Value = "女の子|おんなのこ"
$ws.Cells.Item(93, 1).Value = "boy"
$ws.Cells.Item(93, 2).Value = "男の子|おとこのこ"
$ws.Cells.Item(94, 1).Value = "female student"
$ws.Cells.Item(94, 2).Value = "女子学生|じょしがくせい"
$ws.Cells.Item(95, 1).Value = "small"
$ws.Cells.Item(95, 2).Value = "小さい|ちいさい"
$ws.Cells.Item(96, 1).Value = "elementary school"
$ws.Cells.Item(96, 2).Value = "小学校|しょうがっこう"
$ws.Cells.Item(97, 1).Value = "elementary school student"
$ws.Cells.Item(97, 2).Value = "小学生|しょうがくせい"
$ws.Cells.Item(98, 1).Value = "to meet"
$ws.Cells.Item(98, 2).Value = "会う|あう"
$ws.Cells.Item(99, 1).Value = "company"
$ws.Cells.Item(99, 2).Value = "会社|かいしゃ"
$ws.Cells.Item(100, 1).Value = "office worker"
$ws.Cells.Item(100, 2).Value = "会社員|かいしゃいん"
$ws.Cells.Item(101, 1).Value = "meeting"
$ws.Cells.Item(101, 2).Value = "会議|かいぎ"
$ws.Cells.Item(102, 1).Value = "church"
$ws.Cells.Item(102, 2).Value = "教会|きょうかい"
$ws.Cells.Item(103, 1).Value = "shrine"
$ws.Cells.Item(103, 2).Value = "神社|じんじゃ"
$ws.Cells.Item(104, 1).Value = "society"
$ws.Cells.Item(104, 2).Value = "社会|しゃかい"
$ws.Cells.Item(105, 1).Value = "entry to a company"
$ws.Cells.Item(105, 2).Value = "入社|にゅうしゃ"
$ws.Cells.Item(106, 1).Value = "(my) father"
$ws.Cells.Item(106, 2).Value = "父|ちち"
$ws.Cells.Item(107, 1).Value = "father"
$ws.Cells.Item(107, 2).Value = "お父さん|おとうさん"
$ws.Cells.Item(108, 1).Value = "father and mother"
$ws.Cells.Item(108, 2).Value = "父母|ふぼ"
$ws.Cells.Item(109, 1).Value = "grandfather"
$ws.Cells.Item(109, 2).Value = "祖父|そふ"
$ws.Cells.Item(110, 1).Value = "(my) mother"
$ws.Cells.Item(110, 2).Value = "母|はは"
$ws.Cells.Item(111, 1).Value = "mother"
$ws.Cells.Item(111, 2).Value = "お母さん|おかあさん"
$ws.Cells.Item(112, 1).Value = "mother tongue"
$ws.Cells.Item(112, 2).Value = "母語|ぼご"
$ws.Cells.Item(113, 1).Value = "grandmother"
$ws.Cells.Item(113, 2).Value = "祖母|そぼ"
$ws.Cells.Item(114, 1).Value = "expensive; high"
$ws.Cells.Item(114, 2).Value = "高い|たかい"
$ws.Cells.Item(115, 1).Value = "high school"
$ws.Cells.Item(115, 2).Value = "高校|こうこう"
$ws.Cells.Item(116, 1).Value = "high school student"
$ws.Cells.Item(116, 2).Value = "高校生|こうこうせい"
$ws.Cells.Item(117, 1).Value = "the best"
$ws.Cells.Item(117, 2).Value = "最高|さいこう"
$ws.Cells.Item(118, 1).Value = "school"
$ws.Cells.Item(118, 2).Value = "学校|がっこう"
$ws.Cells.Item(119, 1).Value = "junior high school"
$ws.Cells.Item(119, 2).Value = "中学校|ちゅうがっこう"
$ws.Cells.Item(120, 1).Value = "every day"
$ws.Cells.Item(120, 2).Value = "毎日|まいにち"
$ws.Cells.Item(121, 1).Value = "every week"
$ws.Cells.Item(121, 2).Value = "毎週|まいしゅう"
$ws.Cells.Item(122, 1).Value = "every night"
$ws.Cells.Item(122, 2).Value = "毎晩|まいばん"
$ws.Cells.Item(123, 1).Value = "every year"
$ws.Cells.Item(123, 2).Value = "毎年|まいねん／まいとし"
$ws.Cells.Item(124, 1).Value = "Japanese language"
$ws.Cells.Item(124, 2).Value = "日本語|にほんご"
$ws.Cells.Item(125, 1).Value = "English language"
$ws.Cells.Item(125, 2).Value = "英語|えいご"
$ws.Cells.Item(126, 1).Value = "honorific expressions"
$ws.Cells.Item(126, 2).Value = "敬語|けいご"
$ws.Cells.Item(127, 1).Value = "literature"
$ws.Cells.Item(127, 2).Value = "文学|ぶんがく"
$ws.Cells.Item(128, 1).Value = "composition"
$ws.Cells.Item(128, 2).Value = "作文|さくぶん"
$ws.Cells.Item(129, 1).Value = "letter; character"
$ws.Cells.Item(129, 2).Value = "文字|もじ"
$ws.Cells.Item(130, 1).Value = "culture"
$ws.Cells.Item(130, 2).Value = "文化|ぶんか"
$ws.Cells.Item(131, 1).Value = "grammar"
$ws.Cells.Item(131, 2).Value = "文法|ぶんぽう"
$ws.Cells.Item(132, 1).Value = "to return"
$ws.Cells.Item(132, 2).Value = "帰る|かえる"
$ws.Cells.Item(133, 1).Value = "going home"
$ws.Cells.Item(133, 2).Value = "帰国|きこく"
$ws.Cells.Item(134, 1).Value = "returning home"
$ws.Cells.Item(134, 2).Value = "帰宅|きたく"
$ws.Cells.Item(135, 1).Value = "return"
$ws.Cells.Item(135, 2).Value = "帰り|かえり"
$ws.Cells.Item(136, 1).Value = "to enter"
$ws.Cells.Item(136, 2).Value = "入る|はいる"
$ws.Cells.Item(137, 1).Value = "entrance"
$ws.Cells.Item(137, 2).Value = "入口|いりぐち"
$ws.Cells.Item(138, 1).Value = "to put something in"
$ws.Cells.Item(138, 2).Value = "入れる|いれる"
$ws.Cells.Item(139, 1).Value = "import"
$ws.Cells.Item(139, 2).Value = "輸入|ゆにゅう"
